$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): rotate text 45 degrees and grow row height ---
$ws.Range("A1:G1").Orientation = 45
$ws.Rows.Item(1).RowHeight = 69

# --- Column widths ---
# Target stored widths: A=14.5, B=32.33203125, C=16.33203125, D:G=14.5 (H stays 18, unchanged)
# The engine's ColumnWidth setter quantizes to 1/6-character steps with a
# +5/6 padding offset, so we pre-compensate to land as close as possible
# to the exact target stored widths.
$ws.Columns.Item(1).ColumnWidth = 13.666666666666666
$ws.Columns.Item(2).ColumnWidth = 31.5
$ws.Columns.Item(3).ColumnWidth = 15.5
$ws.Range("D1:G1").ColumnWidth = 13.666666666666666

# --- New "Eigen risico" block (mirrors the existing Bonus/malus block) ---
$ws.Range("A6").Value = "Eigen risico"
$ws.Range("B6").Value = "€ 10043  10044"

$ws.Range("B7").Value = 10043
$ws.Range("D7").Value = "rechts"
$ws.Range("E7").Value = "verwijderen"

$ws.Range("B8").Value = 10044
$ws.Range("D8").Value = "links"
$ws.Range("E8").Value = "niet verwijderen"

# --- New "n/a" block ---
$ws.Range("A9").Value = "n/a"
$ws.Range("B9").Value = "10631 10632 10630          10633"

$ws.Range("C10").Value = 10631
$ws.Range("D10").Value = "links"
$ws.Range("E10").Value = "verwijderen"

$ws.Range("C11").Value = 10632
$ws.Range("D11").Value = "links"
$ws.Range("E11").Value = "verwijderen"

$ws.Range("C12").Value = 10630
$ws.Range("D12").Value = "links"
$ws.Range("E12").Value = "verwijderen"

$ws.Range("C13").Value = 10633
$ws.Range("D13").Value = "links"
$ws.Range("E13").Value = "verwijderen"
